$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:F (mean, lower_80, lower_95, upper_80, upper_95), rows 2-36
$bcdef = New-Object 'object[,]' 35,5
$bcdef[0,0] = 669.123021601995
$bcdef[0,1] = -2545.95428984576
$bcdef[0,2] = -4246.77033123866
$bcdef[0,3] = 3690.31222591602
$bcdef[0,4] = 5952.93766605698
$bcdef[1,0] = 721.464462602777
$bcdef[1,1] = -3735.82446594041
$bcdef[1,2] = -5635.54471469925
$bcdef[1,3] = 5164.50402721429
$bcdef[1,4] = 7128.21459790502
$bcdef[2,0] = 805.465653200318
$bcdef[2,1] = -4385.7965138507
$bcdef[2,2] = -7219.11332544199
$bcdef[2,3] = 6244.19390624855
$bcdef[2,4] = 8911.55956208006
$bcdef[3,0] = 716.820880584684
$bcdef[3,1] = -5266.17793345565
$bcdef[3,2] = -8116.64133130306
$bcdef[3,3] = 7058.37116008172
$bcdef[3,4] = 9534.18110903784
$bcdef[4,0] = 877.765993252498
$bcdef[4,1] = -5736.33179020818
$bcdef[4,2] = -8972.21045276221
$bcdef[4,3] = 7225.71302887793
$bcdef[4,4] = 10541.0210352846
$bcdef[5,0] = 697.192732762604
$bcdef[5,1] = -6369.20931216214
$bcdef[5,2] = -8841.05459453352
$bcdef[5,3] = 7538.17990629029
$bcdef[5,4] = 11460.3467252102
$bcdef[6,0] = 731.902539784532
$bcdef[6,1] = -6872.37854407893
$bcdef[6,2] = -10906.3876887174
$bcdef[6,3] = 7696.28412892197
$bcdef[6,4] = 11421.2421419675
$bcdef[7,0] = 1239.33974742786
$bcdef[7,1] = -6808.3942477251
$bcdef[7,2] = -10213.3412840834
$bcdef[7,3] = 9647.13866708382
$bcdef[7,4] = 12721.5170141943
$bcdef[8,0] = 3291.87492072058
$bcdef[8,1] = -5484.97340230747
$bcdef[8,2] = -9019.18288474367
$bcdef[8,3] = 12363.8026402198
$bcdef[8,4] = 15561.2752436982
$bcdef[9,0] = 3436.8321831056
$bcdef[9,1] = -4068.39438485392
$bcdef[9,2] = -10039.11334277
$bcdef[9,3] = 12558.4810830386
$bcdef[9,4] = 17529.8569962622
$bcdef[10,0] = 1049.93797984718
$bcdef[10,1] = -8097.58454051477
$bcdef[10,2] = -13879.2534491069
$bcdef[10,3] = 10612.0224710846
$bcdef[10,4] = 15727.5999318542
$bcdef[11,0] = 714.448852933579
$bcdef[11,1] = -9114.74884046436
$bcdef[11,2] = -14828.3336341997
$bcdef[11,3] = 10800.6683068343
$bcdef[11,4] = 15459.2931587364
$bcdef[12,0] = 708.358492841157
$bcdef[12,1] = -10374.6536733719
$bcdef[12,2] = -14477.2127560174
$bcdef[12,3] = 11193.3760567504
$bcdef[12,4] = 16067.8195040805
$bcdef[13,0] = 720.476964133274
$bcdef[13,1] = -10374.3975709764
$bcdef[13,2] = -15960.5158899551
$bcdef[13,3] = 12195.7781252619
$bcdef[13,4] = 17107.5297270492
$bcdef[14,0] = 737.82817378964
$bcdef[14,1] = -10912.1933440283
$bcdef[14,2] = -17269.685790999
$bcdef[14,3] = 12902.4053898598
$bcdef[14,4] = 17443.0501711101
$bcdef[15,0] = 961.828990504403
$bcdef[15,1] = -11433.0364749912
$bcdef[15,2] = -18379.9884524787
$bcdef[15,3] = 12716.2491407895
$bcdef[15,4] = 19587.0888400727
$bcdef[16,0] = 927.067360477741
$bcdef[16,1] = -12237.4954051026
$bcdef[16,2] = -19555.8911361781
$bcdef[16,3] = 13418.8130915111
$bcdef[16,4] = 19833.2547737286
$bcdef[17,0] = 766.820368367952
$bcdef[17,1] = -12366.8429284656
$bcdef[17,2] = -19265.1082597988
$bcdef[17,3] = 13823.4740404374
$bcdef[17,4] = 18689.7768003389
$bcdef[18,0] = 755.271208881527
$bcdef[18,1] = -12162.2556994908
$bcdef[18,2] = -19411.2992455641
$bcdef[18,3] = 13293.4877588632
$bcdef[18,4] = 18148.5517611925
$bcdef[19,0] = 1241.72502296714
$bcdef[19,1] = -12522.1695711801
$bcdef[19,2] = -20290.6562565036
$bcdef[19,3] = 14065.5662109771
$bcdef[19,4] = 20163.7404204649
$bcdef[20,0] = 3209.20649013765
$bcdef[20,1] = -9953.8080487423
$bcdef[20,2] = -19304.3341532947
$bcdef[20,3] = 17239.494986226
$bcdef[20,4] = 23824.037641671
$bcdef[21,0] = 3381.19240877192
$bcdef[21,1] = -10297.6756388554
$bcdef[21,2] = -18768.5526540388
$bcdef[21,3] = 17332.2211557377
$bcdef[21,4] = 23707.9458197364
$bcdef[22,0] = 550.159128203021
$bcdef[22,1] = -13656.5421272132
$bcdef[22,2] = -20452.7847027559
$bcdef[22,3] = 14693.1179033266
$bcdef[22,4] = 21972.8037871084
$bcdef[23,0] = 491.855489026024
$bcdef[23,1] = -14563.957693186
$bcdef[23,2] = -22598.6623068145
$bcdef[23,3] = 14686.9712964171
$bcdef[23,4] = 21275.9342061248
$bcdef[24,0] = 550.345039655225
$bcdef[24,1] = -14224.6436162972
$bcdef[24,2] = -21397.4971712646
$bcdef[24,3] = 15056.6652963193
$bcdef[24,4] = 21166.9361822862
$bcdef[25,0] = 612.289933272974
$bcdef[25,1] = -14574.0183375087
$bcdef[25,2] = -22586.8096354868
$bcdef[25,3] = 15587.1423846824
$bcdef[25,4] = 23411.9743159825
$bcdef[26,0] = 628.142908199534
$bcdef[26,1] = -14510.8807863536
$bcdef[26,2] = -22359.5399768796
$bcdef[26,3] = 15692.6037509214
$bcdef[26,4] = 25466.1293186131
$bcdef[27,0] = 575.050931699212
$bcdef[27,1] = -14343.5842594138
$bcdef[27,2] = -23337.9856533679
$bcdef[27,3] = 16808.5686156551
$bcdef[27,4] = 24890.0261684977
$bcdef[28,0] = 488.765499368077
$bcdef[28,1] = -13788.5053881824
$bcdef[28,2] = -24437.6212568924
$bcdef[28,3] = 17614.3720409476
$bcdef[28,4] = 24781.363596329
$bcdef[29,0] = 571.969455280402
$bcdef[29,1] = -14359.0622595462
$bcdef[29,2] = -24015.8795627715
$bcdef[29,3] = 18386.9102087228
$bcdef[29,4] = 25452.5539088846
$bcdef[30,0] = 723.760163076282
$bcdef[30,1] = -13769.7906059638
$bcdef[30,2] = -24857.573729403
$bcdef[30,3] = 18870.1774146288
$bcdef[30,4] = 27166.0586503967
$bcdef[31,0] = 1381.57474398904
$bcdef[31,1] = -14372.56004875
$bcdef[31,2] = -22454.8950526079
$bcdef[31,3] = 19431.4819322545
$bcdef[31,4] = 25806.3303690091
$bcdef[32,0] = 3532.99383452576
$bcdef[32,1] = -12781.7031819401
$bcdef[32,2] = -20662.2473830704
$bcdef[32,3] = 21399.8983673787
$bcdef[32,4] = 28400.9083486306
$bcdef[33,0] = 3807.88087533813
$bcdef[33,1] = -12865.9195145315
$bcdef[33,2] = -21432.9736883286
$bcdef[33,3] = 22092.7607746959
$bcdef[33,4] = 28712.4118274571
$bcdef[34,0] = 1214.92840882248
$bcdef[34,1] = -14916.0871065994
$bcdef[34,2] = -23148.4088217465
$bcdef[34,3] = 19079.1018098617
$bcdef[34,4] = 28842.4538163386

$ws.Range("B2:F36").Value = $bcdef

# New values for column I (diff), rows 2-36
$icol = New-Object 'object[,]' 35,1
$icol[0,0] = 598.123021601995
$icol[1,0] = 704.464462602777
$icol[2,0] = 789.465653200318
$icol[3,0] = 710.820880584684
$icol[4,0] = 873.765993252498
$icol[5,0] = 693.192732762604
$icol[6,0] = 676.902539784532
$icol[7,0] = 1216.33974742786
$icol[8,0] = 3044.87492072058
$icol[9,0] = 3149.8321831056
$icol[10,0] = 986.937979847179
$icol[11,0] = 705.448852933579
$icol[12,0] = 703.358492841157
$icol[13,0] = 718.476964133274
$icol[14,0] = 737.82817378964
$icol[15,0] = 957.828990504403
$icol[16,0] = 922.067360477741
$icol[17,0] = 758.820368367952
$icol[18,0] = 752.271208881527
$icol[19,0] = 1238.72502296714
$icol[20,0] = 3205.20649013765
$icol[21,0] = 3375.19240877192
$icol[22,0] = 545.159128203021
$icol[23,0] = 488.855489026024
$icol[24,0] = 550.345039655225
$icol[25,0] = 612.289933272974
$icol[26,0] = 626.142908199534
$icol[27,0] = 575.050931699212
$icol[28,0] = 486.765499368077
$icol[29,0] = 570.969455280402
$icol[30,0] = 720.760163076282
$icol[31,0] = 1380.57474398904
$icol[32,0] = 3504.99383452576
$icol[33,0] = 3481.88087533813
$icol[34,0] = 1040.92840882248

$ws.Range("I2:I36").Value = $icol

